# Part 4 aliasing bug fix
#
# Slide 16 ("Part 4: Frequency Domain Convolutions") has two placeholder
# text boxes with bracketed student prompts. The second prompt is reworded
# and split into a differently-emphasised run in the middle.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# --- PlaceHolder 3 (right-hand prompt box) -----------------------------
# Old: "[Why does our frequency domain representation of a Gaussian not
#       look like a Gaussian itself? How could we adjust the kernel to
#       make these look more similar?]"
# New: "[Try out some different cutoff values for the 2D Gaussian. " +
#      "What relationship do you notice between the cutoff value and the
#       frequency domain representation? Why is that?" + "]"
$sh3 = $s.Shapes.Item(3)
$tr3 = $sh3.TextFrame.TextRange

# Re-purpose the existing run for the first sentence (keeps its run
# formatting: Arial 14pt, accent2 lumOff, not bold, no strike).
$tr3.Runs(1).Text = "[Try out some different cutoff values for the 2D Gaussian. "

# Append the question as its own run (inherits formatting from the run
# before it: Arial 14pt, accent2 lumOff, not bold, no strike).
[void]$tr3.InsertAfter("What relationship do you notice between the cutoff value and the frequency domain representation? Why is that?")

# Close the bracket in a run matching the original (non-bold) style.
[void]$tr3.InsertAfter("]")
